# Generate Report for Handback
# Adds a second handed-back file (GUID2) alongside the refreshed first file
# (GUID1, which replaces the old GUID in-place with new hash/timestamps) to
# the Overview / zh-cn / de-de report sheets, growing every table by one row.

$wb = $excel.ActiveWorkbook

$GUID1 = "70ddb21b-ec54-4103-819b-9f7406e2b035"
$GUID2 = "8bd78fda-8183-49ca-a9bf-bde54cd5efb7"
$HASH1 = "fbf8c160086cf3820f1c35828195f32c51160451"
$HASH2 = "1e64ca1ae72b1ccdedeb9ad57337ae0850f7c636"

function Col-Index($col) {
    return ([int][char]$col) - ([int][char]'A') + 1
}

# Writes $value into row $row / column-letter $col of $ws as literal text
# (leading "'" stops True/False/date-shaped strings turning into bool/number
# cells) -- unless $value is empty, in which case nothing is written (keeps
# parity with the source workbook's sparser rows, e.g. Overview!D).
function Set-Text($ws, $row, $col, $value) {
    if ($value -eq "") { return }
    $ws.Cells.Item($row, (Col-Index $col)).Value = "'" + $value
}

# Same, but always writes (even "") as literal text -- used for the columns
# that are present-but-empty shared strings in the target (L/N/P).
function Set-TextAlways($ws, $row, $col, $value) {
    $ws.Cells.Item($row, (Col-Index $col)).Value = "'" + $value
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item("Overview")

function Fill-OverviewRow($row, $guid, $dateStr) {
    Set-Text $wsOv $row "A" ($guid + ".md")
    Set-Text $wsOv $row "C" ".md"
    Set-Text $wsOv $row "E" "Handed back: in sync with en-US"
    Set-Text $wsOv $row "F" "Handed back: in sync with en-US"
    Set-Text $wsOv $row "G" $dateStr
    $wsOv.Hyperlinks.Add($wsOv.Cells.Item($row, (Col-Index "B")), `
        ("https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cb8f71f7e727a64e003c6754c7d22b192a7a360f/e2e/" + $guid + ".md"), `
        $null, $null, ("e2e\" + $guid + ".md")) | Out-Null
}

# Row 2 already exists -- refresh its text in place (old GUID -> GUID1).
Fill-OverviewRow 2 $GUID1 "2016-08-27 19:08:06"

# Row 3 is new -- expand the table, then populate it the same way.
$loOv.ListRows.Add() | Out-Null
Fill-OverviewRow 3 $GUID2 "2016-08-27 19:08:06"

# ---------------------------------------------------------------------
# zh-cn / de-de sheets (identical column layout; only the xlf suffix,
# the hyperlink repo, and a couple of datetime columns differ)
# ---------------------------------------------------------------------
function Fill-LangRow($ws, $row, $guid, $hash, $lang, $contentDup, $handoffDate, $handbackDate, $hyperlinkRepo, $hyperlinkCommit) {
    Set-Text $ws $row "A" ($guid + ".md")
    Set-Text $ws $row "B" ".md"
    Set-Text $ws $row "C" "Handed back: in sync with en-US"
    Set-Text $ws $row "D" "e2e"
    Set-Text $ws $row "E" "ht"
    Set-Text $ws $row "F" $contentDup
    Set-Text $ws $row "G" ($guid + "." + $hash + "." + $lang + ".xlf")
    Set-Text $ws $row "H" $handoffDate
    Set-Text $ws $row "I" ($guid + ".md")
    Set-Text $ws $row "J" ($guid + "." + $hash + "." + $lang + ".xlf")
    Set-Text $ws $row "K" $handbackDate
    Set-TextAlways $ws $row "L" ""
    Set-Text $ws $row "M" "True"
    Set-TextAlways $ws $row "N" ""
    Set-Text $ws $row "O" "False"
    Set-TextAlways $ws $row "P" ""

    $ws.Hyperlinks.Add($ws.Cells.Item($row, (Col-Index "A")), `
        ("https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cb8f71f7e727a64e003c6754c7d22b192a7a360f/e2e/" + $guid + ".md"), `
        $null, $null, ($guid + ".md")) | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item($row, (Col-Index "I")), `
        ("https://github.com/" + $hyperlinkRepo + "/blob/" + $hyperlinkCommit + "/e2e/" + $guid + ".md"), `
        $null, $null, ($guid + ".md")) | Out-Null
}

# --- zh-cn ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item("zh-cn")

Fill-LangRow $wsZh 2 $GUID1 $HASH1 "zh-cn" "False" "2016-08-27 19:07:58" "2016-08-27 19:08:27" `
    "OpenLocalizationTestOrg/ol-test0-zhcn" "0fc04650e2bfb2d30c4930422273ded30f5b0110"

$loZh.ListRows.Add() | Out-Null
Fill-LangRow $wsZh 3 $GUID2 $HASH2 "zh-cn" "True" "2016-08-27 19:07:58" "2016-08-27 19:08:27" `
    "OpenLocalizationTestOrg/ol-test0-zhcn" "0fc04650e2bfb2d30c4930422273ded30f5b0110"

# --- de-de ---
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item("de-de")

Fill-LangRow $wsDe 2 $GUID1 $HASH1 "de-de" "False" "2016-08-27 19:08:06" "2016-08-27 19:08:34" `
    "OpenLocalizationTestOrg/ol-test0-dede" "878c59f8885a98319c6caab8df8967ba24a5a87b"

$loDe.ListRows.Add() | Out-Null
Fill-LangRow $wsDe 3 $GUID2 $HASH2 "de-de" "True" "2016-08-27 19:08:06" "2016-08-27 19:08:34" `
    "OpenLocalizationTestOrg/ol-test0-dede" "878c59f8885a98319c6caab8df8967ba24a5a87b"

Write-Output "done"
